# DR solution updated in last page
# Updates the "Backend Bandwidth Allocation" table on the last slide (slide 5):
#   - header columns renamed (Backend -> Api Engineer, Frontend -> Key Activities)
#   - column widths re-distributed
#   - new content filled in for the previously-empty rows

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tbl = $sh.Table

# --- Column widths (EMU -> points, 1 pt = 12700 EMU) ---
$tbl.Columns.Item(1).Width = 5500914 / 12700
$tbl.Columns.Item(2).Width = 5152572 / 12700
$tbl.Columns.Item(3).Width = 10814595 / 12700

# --- Header row ---
$tbl.Cell(1, 2).Shape.TextFrame.TextRange.Text = "Api Engineer"
$tbl.Cell(1, 3).Shape.TextFrame.TextRange.Text = "Key Activities"

# --- CMDB row ---
$tbl.Cell(2, 2).Shape.TextFrame.TextRange.Text = "2"
$tbl.Cell(2, 3).Shape.TextFrame.TextRange.Text = "Open Api specification / Data Design +`rSpring Boot Api’s + GitOps "

# --- Security/RBAC row ---
$tbl.Cell(3, 2).Shape.TextFrame.TextRange.Text = "1"
$tbl.Cell(3, 3).Shape.TextFrame.TextRange.Text = "User/Role/Group/Transaction/ Session / Authentication/ MFA / Activity Logs"

# --- Workflow row ---
$tbl.Cell(4, 1).Shape.TextFrame.TextRange.Text = "Workflow"
$tbl.Cell(4, 2).Shape.TextFrame.TextRange.Text = "3"
$tbl.Cell(4, 3).Shape.TextFrame.TextRange.Text = "Automation flows using Aws steps and lambdas (Golang / node / python)"

# --- Monitoring row ---
$tbl.Cell(5, 1).Shape.TextFrame.TextRange.Text = "Monitoring"
$tbl.Cell(5, 2).Shape.TextFrame.TextRange.Text = "2"
$tbl.Cell(5, 3).Shape.TextFrame.TextRange.Text = "Write routines for continuous data collection and implements algorithms deduce SLA’s"

# --- Analytics row ---
$tbl.Cell(6, 1).Shape.TextFrame.TextRange.Text = "Analytics"
$tbl.Cell(6, 2).Shape.TextFrame.TextRange.Text = "2"
$tbl.Cell(6, 3).Shape.TextFrame.TextRange.Text = "Write Reports and Analytics "
